# Add a new "prop" worksheet (ethnicity proportions, derived from the
# "ethnicities" sheet) after "ethnicities", make it the active sheet/tab,
# and update the previously-active "ethnicities" sheet's selection.

$wb = $excel.ActiveWorkbook

# --- Update selection on the "ethnicities" sheet (loses tabSelected once
#     a later sheet is activated below) -------------------------------
$wsEth = $wb.Worksheets.Item("ethnicities")
$wsEth.Range("G28").Select() | Out-Null

# --- Create the new "prop" sheet as the last tab ----------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $last)
$ws.Name = "prop"

# Header row (reuses existing shared strings; note the original author's
# off-by-one header labels vs. the data columns - reproduced verbatim)
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Value = "White"
$ws.Range("C1").Value = "African American"
$ws.Range("D1").Value = "Asian"
$ws.Range("E1").Value = "Other"
$ws.Range("F1").Value = "Unknown"

# Row 2: Cases
$ws.Range("A2").Value = "Cases"
$ws.Range("B2").Value = 47.244934233913973
$ws.Range("C2").Value = 32.91859225026662
$ws.Range("D2").Value = 10.344827586206897
$ws.Range("E2").Value = 2.0618556701030926
$ws.Range("F2").Value = 3.2705296836118025

# Row 3: Hospitalizations
$ws.Range("A3").Value = "Hospitalizations"
$ws.Range("B3").Value = 49.536178107606673
$ws.Range("C3").Value = 35.807050092764378
$ws.Range("D3").Value = 9.833024118738404
$ws.Range("E3").Value = 1.1131725417439702
$ws.Range("F3").Value = 3.1539888682745829

# Row 4: ICU
$ws.Range("A4").Value = "ICU"
$ws.Range("B4").Value = 47.029702970297024
$ws.Range("C4").Value = 39.10891089108911
$ws.Range("D4").Value = 9.4059405940594054
$ws.Range("E4").Value = 1.9801980198019802
$ws.Range("F4").Value = 1.9801980198019802

# Row 5: Deaths
$ws.Range("A5").Value = "Deaths"
$ws.Range("B5").Value = 71.090047393364927
$ws.Range("C5").Value = 18.009478672985782
$ws.Range("D5").Value = 6.6350710900473935
$ws.Range("E5").Value = 1.8957345971563981
$ws.Range("F5").Value = 0.94786729857819907

# Column widths matching the source file (best-fit look: 14 / 15 chars).
# ColumnWidth setter stores a value offset by +5/6 internally, so back it
# out here to land exactly on width="14" / width="15" after save.
$ws.Columns.Item(1).ColumnWidth = 13.166666666666666
$ws.Columns.Item(3).ColumnWidth = 14.166666666666666

# Selection/active cell on the new sheet - selecting it makes "prop" the
# active tab (and clears tabSelected on "ethnicities").
$ws.Range("L9").Select() | Out-Null
